# Apply cryptos list price/volume updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.772.99"
$ws.Range("E2").Value = "  -2.20%  "
$ws.Range("D3").Value = "3.469.67"
$ws.Range("E3").Value = "  -2.64%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.59%  "
$ws.Range("D7").Value = "3.465.81"
$ws.Range("E7").Value = "  -2.72%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.477"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.67%  "
$ws.Range("E10").Value = "  -3.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.55"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000211"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.00%  "
$ws.Range("D14").Value = "4.056.35"
$ws.Range("E14").Value = "  -2.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.21%  "
$ws.Range("D16").Value = "3.483.77"
$ws.Range("E16").Value = "  -2.31%  "
$ws.Range("D17").Value = "66.755.73"
$ws.Range("E17").Value = "  -2.12%  "
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("E19").Value = "  -5.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "431.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.603"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value = "3.608.67"
$ws.Range("E26").Value = "  -2.54%  "
$ws.Range("E27").Value = "  -8.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.26%  "
$ws.Range("E30").Value = "  -3.90%  "
$ws.Range("E31").Value = "  -7.36%  "
$ws.Range("E32").Value = "  -3.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.993"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.18"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.65%  "
$ws.Range("D35").Value = "3.460.45"
$ws.Range("E35").Value = "  -2.69%  "
$ws.Range("E36").Value = "  -8.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.05%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.49%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "174.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0875"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.81%  "
$ws.Range("E43").Value = "  -11.00%  "
$ws.Range("E44").Value = "  -4.00%  "
$ws.Range("E45").Value = "  -0.78%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.15%  "
$ws.Range("E48").Value = "  -8.89%  "
$ws.Range("E49").Value = "  -5.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -9.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.969"
$ws.Range("D51").Style = "Normal"
